# Generate and save output file after processing
# Insert three new columns (history, electives, cs) into the general_college_subjects
# block, right before the "arts" column, and refresh the row-2 text casing / values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "arts" currently lives in column R (18). Insert 3 new columns right there so the
# existing "arts" column (and everything after it) shifts right by 3, becoming U.
$ws.Range("R1:T1").EntireColumn.Insert()

# New header cells for the inserted columns.
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# New data cells (row 2) for the inserted columns.
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 3
$ws.Range("T2").Value = 0

# Normalize casing of the existing text values in row 2.
$ws.Range("D2").Value = "considered"
$ws.Range("E2").Value = "important"
$ws.Range("F2").Value = "considered"
$ws.Range("G2").Value = "considered"
$ws.Range("H2").Value = "very important"
$ws.Range("I2").Value = "very important"
$ws.Range("J2").Value = "considered"

$wb.Save()
